$p = $ppt.ActivePresentation
$p.Slides.Item(2).Delete()
Write-Output $p.Slides.Count
